# Commit: "Added consistency to unit naming for web display (v0.17)"
#
# On the "all-variables" sheet, column K holds the human-readable "units"
# label for each series. A handful of ad-hoc unit strings are standardized
# onto a smaller, consistent vocabulary shared with the rest of the sheet
# (and the shared string pool). Two labels that become fully unused after
# this pass ("$bn" and "$millions") naturally drop out of the workbook the
# same way Excel itself would prune them on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all-variables")
$ws.Activate()

# row -> new unit text, grouped by the standardized label they adopt
$billionsOf2012 = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,109,110)
foreach ($r in $billionsOf2012) {
    $ws.Range("K$r").Value = "billions of 2012 `$ "
}

$thousands = @(51,58)
foreach ($r in $thousands) {
    $ws.Range("K$r").Value = "thousands"
}

$millionsDollars = @(54,56,68,113,114,115,116,117,118,119,120,121,122,123,124,125)
foreach ($r in $millionsDollars) {
    $ws.Range("K$r").Value = "millions `$"
}

$ws.Range("K57").Value = "millions"
$ws.Range("K61").Value = "`$ per barrel"

$billionsDollars = @(62,63,64,112)
foreach ($r in $billionsDollars) {
    $ws.Range("K$r").Value = "billions `$"
}

$ws.Range("K70").Value = "index"
$ws.Range("K72").Value = "index (2012 = 100)"

# View state: the sheet had scrolled down to row 93 with E127 selected;
# the edit leaves it scrolled near the top of the table with K58 selected.
$ws.Range("K58").Select()
